# Daily attendance processing - 2026-02-01 18:08:04
# Normalize the "Recorded By" column (G): when the value is a comma-separated
# list whose first entry is literally "System", move that first entry to the
# end of the list (swap the first and last entries), keeping every entry's
# original text/casing intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }

    $parts = $val -split ','
    if ($parts.Count -lt 2) { continue }

    for ($i = 0; $i -lt $parts.Count; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts[0] -ceq 'System') {
        $first = $parts[0]
        $last = $parts[$parts.Count - 1]
        $parts[0] = $last
        $parts[$parts.Count - 1] = $first
        $cell.Value2 = [string]::Join(', ', $parts)
    }
}
